# Update the date heading
$d = $word.ActiveDocument
$d.Content.Find.Execute("2026-01-08 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-09 Friday", 2) | Out-Null

# Update the uniquely-identifiable multiplication problems throughout the
# document (each search string occurs exactly once).
$d.Content.Find.Execute("303×5=", $true, $false, $false, $false, $false, $true, 1, $false, "731×6=", 2) | Out-Null
$d.Content.Find.Execute("208×9=", $true, $false, $false, $false, $false, $true, 1, $false, "891×6=", 2) | Out-Null
$d.Content.Find.Execute("865×5=", $true, $false, $false, $false, $false, $true, 1, $false, "691×6=", 2) | Out-Null
$d.Content.Find.Execute("145×2=", $true, $false, $false, $false, $false, $true, 1, $false, "271×3=", 2) | Out-Null
$d.Content.Find.Execute("641×2=", $true, $false, $false, $false, $false, $true, 1, $false, "514×9=", 2) | Out-Null

$d.Content.Find.Execute("304×8=", $true, $false, $false, $false, $false, $true, 1, $false, "134×9=", 2) | Out-Null
$d.Content.Find.Execute("864×7=", $true, $false, $false, $false, $false, $true, 1, $false, "694×9=", 2) | Out-Null
$d.Content.Find.Execute("398×4=", $true, $false, $false, $false, $false, $true, 1, $false, "496×5=", 2) | Out-Null
$d.Content.Find.Execute("744×8=", $true, $false, $false, $false, $false, $true, 1, $false, "150×6=", 2) | Out-Null

$d.Content.Find.Execute("402×3=", $true, $false, $false, $false, $false, $true, 1, $false, "915×8=", 2) | Out-Null
$d.Content.Find.Execute("365×3=", $true, $false, $false, $false, $false, $true, 1, $false, "952×5=", 2) | Out-Null
$d.Content.Find.Execute("857×4=", $true, $false, $false, $false, $false, $true, 1, $false, "244×2=", 2) | Out-Null
$d.Content.Find.Execute("218×2=", $true, $false, $false, $false, $false, $true, 1, $false, "182×5=", 2) | Out-Null
$d.Content.Find.Execute("745×3=", $true, $false, $false, $false, $false, $true, 1, $false, "951×6=", 2) | Out-Null

$d.Content.Find.Execute("713×2=", $true, $false, $false, $false, $false, $true, 1, $false, "876×8=", 2) | Out-Null
$d.Content.Find.Execute("852×6=", $true, $false, $false, $false, $false, $true, 1, $false, "299×2=", 2) | Out-Null
$d.Content.Find.Execute("422×8=", $true, $false, $false, $false, $false, $true, 1, $false, "902×6=", 2) | Out-Null
$d.Content.Find.Execute("664×7=", $true, $false, $false, $false, $false, $true, 1, $false, "882×7=", 2) | Out-Null

$d.Content.Find.Execute("736×5=", $true, $false, $false, $false, $false, $true, 1, $false, "752×5=", 2) | Out-Null
$d.Content.Find.Execute("927×4=", $true, $false, $false, $false, $false, $true, 1, $false, "844×6=", 2) | Out-Null
$d.Content.Find.Execute("127×6=", $true, $false, $false, $false, $false, $true, 1, $false, "976×5=", 2) | Out-Null
$d.Content.Find.Execute("992×3=", $true, $false, $false, $false, $false, $true, 1, $false, "410×6=", 2) | Out-Null
$d.Content.Find.Execute("843×2=", $true, $false, $false, $false, $false, $true, 1, $false, "888×3=", 2) | Out-Null

# "456×6=" occurs twice in the table and must map to two different values
# depending on position, so target each occurrence by its specific cell and
# use wdReplaceOne (1) so the replacement cannot bleed into other cells.
$t = $d.Tables.Item(1)
$t.Cell(5, 4).Range.Find.Execute("456×6=", $true, $false, $false, $false, $false, $true, 1, $false, "573×2=", 1) | Out-Null
$t.Cell(15, 2).Range.Find.Execute("456×6=", $true, $false, $false, $false, $false, $true, 1, $false, "550×6=", 1) | Out-Null
